# fm10: Remove R8 from fab files
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rev 5")

# 1. Rename the sheet from "Rev 5" to "Rev 4.1"
$ws.Name = "Rev 4.1"

# 2. Update the revision cell F2 from numeric 5 to text "4.1"
$ws.Range("F2").Value = "4.1"

# 3. Remove "R8" from the Position list in row 18 (item 13, R_0805 jumpers)
$ws.Range("E18").Value = "R2 R10 R31 R55 R56"

# Mark the updated cell red, to indicate removed component (per legend)
$ws.Range("E18").Interior.Color = 3355647

# 4. Restore view state
$ws.Range("A1").Select()
$ws.Range("F3").Select()
